$d = $word.ActiveDocument

$replacements = @(
    @{old="2025-07-18 Friday"; new="2025-07-19 Saturday"},
    @{old="37×47="; new="36×86="},
    @{old="69×37="; new="33×49="},
    @{old="46×29="; new="91×30="},
    @{old="76×30="; new="62×96="},
    @{old="71×81="; new="66×45="},
    @{old="46×65="; new="64×92="},
    @{old="58×77="; new="12×88="},
    @{old="95×22="; new="37×22="},
    @{old="35×41="; new="14×55="},
    @{old="90×78="; new="50×33="},
    @{old="12×94="; new="72×56="},
    @{old="99×76="; new="95×12="},
    @{old="26×34="; new="96×39="},
    @{old="53×56="; new="77×31="},
    @{old="26×73="; new="26×99="},
    @{old="39×85="; new="59×62="},
    @{old="38×26="; new="84×71="},
    @{old="15×47="; new="76×77="},
    @{old="72×93="; new="77×16="},
    @{old="14×35="; new="75×24="},
    @{old="75×41="; new="95×15="},
    @{old="76×53="; new="56×37="},
    @{old="45×50="; new="81×25="},
    @{old="40×55="; new="31×66="},
    @{old="13×98="; new="47×49="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
